# Assignment 3.docx edit
#
# 1) Insert an empty "_GoBack" bookmark (bookmarkStart/bookmarkEnd pair,
#    id 0) at the very start of the document content, right before the
#    title run - this is the bookmark Word silently drops at the last
#    place the document was edited.
# 2) Merge the two adjacent, identically formatted runs that spell out
#    "LED BLINKING CODE" + ":" into a single run reading
#    "LED BLINKING CODE:".

$d = $word.ActiveDocument

# --- 1) "_GoBack" bookmark at the very start of the document --------------
# Word's Bookmarks.Add needs a non-degenerate anchor to reliably place a
# zero-length bookmark at the absolute start of the story, so nudge a
# placeholder character in, bookmark right after it, then remove the
# placeholder again - leaving a clean, empty bookmark at position 0.
$startRng = $d.Range(0, 0)
$startRng.InsertBefore("X")

$bookmarkRng = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRng)

$placeholderRng = $d.Range(0, 1)
$placeholderRng.Delete()

# --- 2) Merge "LED BLINKING CODE" + ":" into a single run ------------------
$d.Content.Find.Execute("LED BLINKING CODE:", $true, $false, $false, $false,
                         $false, $true, 1, $false, "LED BLINKING CODE:", 2)

Write-Output "done"
